$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "k" column (B) is being removed, shifting every column after it one to the left
# and the remaining headers are being re-labelled / reordered.
$optVal  = $ws.Range("C1").Value2   # "Opt"
$heurVal = $ws.Range("D1").Value2   # "Heur"
$gapVal  = $ws.Range("G1").Value2   # "GAP LB-UB"

$ws.Range("B1").Value = $optVal
$ws.Range("C1").Value = $heurVal
$ws.Range("D1").Value = "Best UB"
$ws.Range("E1").Value = "Best LB"
$ws.Range("F1").Value = $gapVal
$ws.Range("G1").Value = "Num it"
$ws.Range("H1").Value = "GAP UB opt"
$ws.Range("I1").Value = "GAP LB opt"
$ws.Range("J1").Value = "GAP uB heut"
$ws.Range("K1").Value = "GAP LB heu"

# Drop the now-unused last column (previously "It")
$ws.Range("L1").ClearContents()

$ws.Range("A1").Select()
